$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: merge "ISYS3001" + " <en-dash> Assignment 1 exercise"
#    into a single run (keeps the existing bold/size formatting).
# ---------------------------------------------------------------------------
$enDash = [char]0x2013
$oldTitle = "ISYS3001 " + $enDash + " Assignment 1 exercise"
$d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $oldTitle, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Remember that this is a public repository - your changes could be
#    seen by anyone who looks!" - merge runs into one.
# ---------------------------------------------------------------------------
$old2 = "Remember that this is a public repository - your changes could be seen by anyone who looks!"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Add some comments about Version management outside this border, or
#    just add some text so there is a change to this file." - merge runs.
# ---------------------------------------------------------------------------
$old3 = "Add some comments about Version management outside this border, or just add some text so there is a change to this file."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Remember that your GitHub user ID must be submitted in your
#    assignment report!" - merge runs.
# ---------------------------------------------------------------------------
$old4 = "Remember that your GitHub user ID must be submitted in your assignment report!"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Final paragraph: "<hellip> This is austin" -> split into a
#    "<hellip>" paragraph and a "This is AustinFYX." paragraph, dropping
#    the lone joining space and the stray language overrides, and drop
#    the empty trailing paragraph after it.
# ---------------------------------------------------------------------------
$ellipsis = [char]0x2026
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($ellipsis)) {
        $targetIdx = $i
    }
}
$p6 = $d.Paragraphs.Item($targetIdx)
$p7 = $d.Paragraphs.Item($targetIdx + 1)
$wholeRange = $d.Range($p6.Range.Start, $p7.Range.End)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$xml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>' + $ellipsis + '</w:t></w:r></w:p>' + `
       '<w:p xmlns:w="' + $wNs + '"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>T</w:t></w:r><w:r><w:t>his is AustinFYX.</w:t></w:r></w:p>'
$wholeRange.InsertXML($xml) | Out-Null
